# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N ("Late") on the
# "Repayment Schedule" sheet to make room for a Variable Instalments
# column, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new column at position N (14), shifting the existing
# "Late" and "Outstanding" columns one place to the right.
$ws.Columns.Item(14).Insert()

# Give the newly inserted column its own width (stored column width of 10).
$ws.Columns.Item(14).ColumnWidth = 9.166666666666666

# Update the saved selection to match the post-edit state.
$ws.Range("T8").Select()
